$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title line: split the bold "First Sprint" run into "First" + " Sprint"
#    and drop a _GoBack bookmark right before it (after "Sprint Backlog: ").
#    "Sprint Backlog: " is always 16 characters long at the very start of
#    the document, so the boundary between the two runs is a fixed offset.
# ---------------------------------------------------------------------------
$boundary = 16  # length of "Sprint Backlog: "

$spaceSprintRange = $d.Range($boundary + 5, $boundary + 12)   # " Sprint"
$spaceSprintRange.Text = ""
$insertPoint = $d.Range($boundary + 5, $boundary + 5)
$insertPoint.InsertAfter(" Sprint")
$newRunRange = $d.Range($boundary + 5, $boundary + 12)
$newRunRange.Font.Bold = $true
$newRunRange.Font.Size = 14

$d.Bookmarks.Add("_GoBack", $d.Range($boundary, $boundary))

# ---------------------------------------------------------------------------
# 2. User Story 1 text + its task rewordings.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("As a User, I want to make a profile so that, my information will display.", $true, $false, $false, $false, $false, $true, 1, $false, "As a User, I want to enter the website so that, I can view homepage and signup there.", 2)

$d.Content.Find.Execute("Design user interface for user profile page.", $true, $false, $false, $false, $false, $true, 1, $false, "Design user interface for user home page & signup form.", 2)

$d.Content.Find.Execute("Write code for user profile page.", $true, $false, $false, $false, $false, $true, 1, $false, "Write code for home page and form.", 2)

$d.Content.Find.Execute("Connect user profile page with database.", $true, $false, $false, $false, $false, $true, 1, $false, "Connect home page with database.", 2)

$d.Content.Find.Execute("Test the user profile page. ", $true, $false, $false, $false, $false, $true, 1, $false, "Test the home page & signup ", 2)

# ---------------------------------------------------------------------------
# 3. Numeric day-count cells that changed value.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)

# Row 1 ("make a profile" story) Day 2 column: 0 -> 1
$day2Cell = $t.Cell(2, 5)
$day2Para = $day2Cell.Range.Paragraphs(1)
$day2Run = $d.Range($day2Para.Range.Start, $day2Para.Range.Start + 1)
$day2Run.Text = "1"

# "Write code" task row Day 3 column: 2 -> 1
$day3Cell = $t.Cell(3, 6)
$day3Para = $day3Cell.Range.Paragraphs(1)
$day3Run = $d.Range($day3Para.Range.Start, $day3Para.Range.Start + 1)
$day3Run.Text = "1"

# ---------------------------------------------------------------------------
# 4. Remove the entire "edit my information" user story (its 4 rows: the
#    story row + 3 follow-on task rows) -- they are the last 4 rows in the
#    table.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
    $t.Rows.Last.Delete()
}

# ---------------------------------------------------------------------------
# 5. The now-last row ("Test the user profile page.") shrinks its height.
# ---------------------------------------------------------------------------
$lastRow = $t.Rows.Last
$lastRow.Height = 34.6   # 692 dxa (twentieths of a point)

# ---------------------------------------------------------------------------
# 6. Two extra trailing empty paragraphs after the table.
# ---------------------------------------------------------------------------
$tableEnd = $t.Range.End
$afterTable = $d.Range($tableEnd, $tableEnd)
$afterTable.InsertParagraphAfter()
$afterTable.InsertParagraphAfter()
